$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-13 (NOMBRE, APELLIDO, NUMERO, ULTIMO_PAGO)
$data = @(
    @("Daniela",   "Villamizar",   3218490916, 44604),
    @("Juan",      "Botero",       3218748814, 44605),
    @("Cristian",  "Solarte",      3148227994, 44606),
    @("Julian",    "Aristizabal",  3046145922, 44607),
    @("Juan",      "Londoño",      3163610054, 44608),
    @("Mauricio",  "Herrera",      3117754781, 44609),
    @("Esteban",   "Meneses",      3108017554, 44610),
    @("Carlos",    "Paraco",       3188288098, 44611),
    @("Alejandra", "Ruiz",         3167357054, 44612),
    @("Daniela",   "Bustos",       3164224295, 44613),
    @("Laura",     "Hoyos",        3113829197, 44614),
    @("Manuela",   "Rojas",        3114244572, 44615)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    if ($row -eq 3) {
        # matches the original authoring order, which registered the
        # surname string before the first name string for this row
        $ws.Cells.Item($row, 2).Value = $vals[1]
        $ws.Cells.Item($row, 1).Value = $vals[0]
    } else {
        $ws.Cells.Item($row, 1).Value = $vals[0]
        $ws.Cells.Item($row, 2).Value = $vals[1]
    }
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
}

# Remove the now-unused trailing rows (previously rows 14 and 15)
$ws.Rows("14:15").Delete()

# Restore selection to match the saved workbook state
[void]$ws.Range("D16").Select()
